$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column "想去人数" counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 752
$ws1.Range("F4").Value = 24
$ws1.Range("F6").Value = 1182
$ws1.Range("F10").Value = 585
$ws1.Range("F14").Value = 97
$ws1.Range("F15").Value = 6
$ws1.Range("F16").Value = 85
$ws1.Range("F17").Value = 290
$ws1.Range("F21").Value = 5981
$ws1.Range("F22").Value = 5294

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 89

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 752
$ws4.Range("F4").Value = 24
$ws4.Range("F6").Value = 1182
$ws4.Range("F10").Value = 585
$ws4.Range("F13").Value = 89
$ws4.Range("F16").Value = 97
$ws4.Range("F17").Value = 6
$ws4.Range("F18").Value = 85
$ws4.Range("F19").Value = 290
$ws4.Range("F23").Value = 5981
$ws4.Range("F25").Value = 5294
